$wb = $excel.ActiveWorkbook

# Data describing each cell update: SheetName -> list of (CellRef, NewValue)
$updates = @{
    "Citywide Totals" = @{ "K2"=5711; "K3"=5861; "D4"=1977; "G4"=1489; "K4"=1220; "K5"=416; "K6"=6471; "D7"=28167; "G7"=24716; "K7"=19679 }
    "Logan Square" = @{ "K3"=63; "K7"=246 }
    "Austin" = @{ "K2"=359; "K3"=398; "K6"=436; "K7"=1303 }
    "Garfield Park" = @{ "K2"=228; "K3"=310; "K7"=846 }
    "West Pullman" = @{ "K2"=113; "K7"=332 }
    "Grand Crossing" = @{ "K2"=194; "K6"=191 }
    "New City" = @{ "K2"=148; "K7"=452 }
    "By Neighborhood" = @{ "K2"=175; "K4"=74; "K8"=1303; "K9"=85; "K11"=376; "K19"=574; "K20"=466; "K21"=65; "K27"=181; "K29"=1067; "K31"=217; "K33"=846; "K36"=255; "K40"=46; "K42"=726; "K43"=169; "K44"=168; "K48"=250; "K51"=251; "K52"=511; "K53"=246; "K54"=379; "K55"=218; "D63"=356; "G63"=290; "K65"=452; "K67"=764; "K68"=51; "K74"=20; "K75"=65; "K76"=269; "K77"=137; "K78"=224; "K79"=494; "K80"=69; "K85"=926; "K86"=127; "K93"=74; "K94"=263; "K95"=332; "K96"=207; "D101"=28167; "G101"=24716; "K101"=19679 }
    "Gage Park" = @{ "K3"=56; "K6"=76; "K7"=217 }
    "North Lawndale" = @{ "K2"=217; "K3"=274; "K4"=43; "K6"=213; "K7"=764 }
    "Loop" = @{ "K3"=94; "K5"=3; "K7"=379 }
    "Englewood" = @{ "K3"=381; "K6"=307; "K7"=1067 }
    "Lake View" = @{ "K2"=36; "K7"=250 }
    "Chatham" = @{ "K3"=177; "K7"=574 }
    "Irving Park" = @{ "K2"=42; "K3"=45; "K7"=168 }
    "River North" = @{ "K2"=60; "K7"=269 }
    "Humboldt Park" = @{ "K2"=194; "K3"=224; "K7"=726 }
    "Rogers Park" = @{ "K2"=67; "K3"=53; "K7"=224 }
    "Lower West Side" = @{ "K6"=74; "K7"=218 }
    "West Ridge" = @{ "K2"=64; "K7"=207 }
    "Chinatown" = @{ "K4"=1; "K7"=65 }
    "Roseland" = @{ "K3"=160; "K7"=494 }
    "Chicago Lawn" = @{ "K2"=154; "K3"=152; "K7"=466 }
    "Grand Boulevard" = @{ "K3"=74; "K6"=58; "K7"=255 }
    "West Lawn" = @{ "K2"=25; "K3"=17; "K7"=74 }
    "West Loop" = @{ "K3"=54; "K7"=263 }
    "Belmont Cragin" = @{ "K2"=132; "K4"=23; "K7"=376 }
    "Avalon Park" = @{ "K6"=22; "K7"=85 }
    "Albany Park" = @{ "K3"=42; "K7"=175 }
    "Edgewater" = @{ "K3"=45; "K7"=181 }
    "Streeterville" = @{ "K4"=54; "K7"=127 }
    "Pullman" = @{ "K6"=12; "K7"=65 }
    "Little Italy, UIC" = @{ "K2"=71; "K4"=28; "K7"=251 }
    "North Park" = @{ "K2"=23; "K7"=51 }
    "Hyde Park" = @{ "K6"=68; "K7"=169 }
    "South Shore" = @{ "K2"=302; "K3"=315; "K7"=926 }
    "Riverdale" = @{ "K3"=55; "K4"=8; "K7"=137 }
    "Rush & Division" = @{ "K6"=34; "K7"=69 }
    "Hegewisch" = @{ "K6"=8; "K7"=46 }
    "Little Village" = @{ "K2"=137; "K3"=147; "K6"=182; "K7"=511 }
    "Archer Heights" = @{ "K3"=15; "K7"=74 }
    "Printers Row" = @{ "K2"=3; "K7"=20 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($ref in $cellMap.Keys) {
        $ws.Range($ref).Value = $cellMap[$ref]
    }
}